# "Drop in results from RMI script"
# Overwrite a handful of outlier cells on the passenger (psgr) sheet so the
# row is uniform across vehicle types/years, matching what the RMI script
# produced, and remove the leftover helper formula in B2.

$wb = $excel.ActiveWorkbook

$wsPsgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

# Row 2 (LDVs): B2 had a leftover "=0.076+(0.076-0.0725)" formula and D2 was
# an outlier (7.35%); both now become a plain 7.6% value like the rest of
# the row.
$wsPsgr.Range("B2").Value = 0.076
$wsPsgr.Range("D2").Value = 0.076

# Row 5 (ships): B5 and E5 were outliers left at 1%; bring them in line with
# the rest of the row at 2.9%.
$wsPsgr.Range("B5").Value = 0.029
$wsPsgr.Range("E5").Value = 0.029

# The script's last action leaves the "About" sheet active/selected instead
# of the passenger sheet.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
